$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Web")

# New "Upload and Download" section (rows 62-67), mirroring the existing
# section pattern (e.g. brokenLinks at rows 58-61).
$ws.Range("A62").Value = "uploadDownload"

$ws.Range("A63").Value = "uploadDownloadNav"
$ws.Range("B63").Value = '//*[@id="item-7"]'
$ws.Range("C63").Value = "By.xpath"

$ws.Range("A64").Value = "uploadDownloadScroll"
$ws.Range("B64").Value = '//*[@id="app"]/div/div/div/div[2]/div[2]/h1'
$ws.Range("C64").Value = "By.xpath"

$ws.Range("A65").Value = "imgDownloadBtn"
$ws.Range("B65").Value = '//*[@id="downloadButton"]'
$ws.Range("C65").Value = "By.xpath"

$ws.Range("A66").Value = "uploadFileInput"
$ws.Range("B66").Value = '//*[@id="app"]/div/div/div/div[2]/div[2]/div[2]/form/div'
$ws.Range("C66").Value = "By.xpath"

$ws.Range("A67").Value = "uploadedFileName"
$ws.Range("B67").Value = '//*[@id="uploadedFilePath"]'
$ws.Range("C67").Value = "By.xpath"

# Update the saved view state (scroll position / selection).
$ws.Application.ActiveWindow.ScrollRow = 55
$ws.Range("B66").Select()
